# Insert new bullet paragraphs describing the boundary estimation algorithm
# right after the "Full-Stack Development and Data Engineering" sub-heading
# paragraph under the Siege Analytics / PARTNER role.

$d = $word.ActiveDocument

# Locate the 1-based index of the anchor paragraph.
$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Trim() -eq "Full-Stack Development and Data Engineering") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Anchor paragraph 'Full-Stack Development and Data Engineering' not found"
}

$newBullets = @(
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times",
    "• Architected systems supporting 2,500+ concurrent users conducting redistricting analysis",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

$insertAfter = $anchorIndex
foreach ($bulletText in $newBullets) {
    $anchorPara = $d.Paragraphs.Item($insertAfter)
    $anchorPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($insertAfter + 1)
    $newPara.Range.Text = $bulletText
    $insertAfter = $insertAfter + 1
}
